$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.341.26'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').Value = '1.876.23'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7118'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.35'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3111'
$ws.Range('E8').Value = '  +1.20%  '
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08466'
$ws.Range('E11').Value = '  +2.58%  '
$ws.Range('D12').Value = '1.914.06'
$ws.Range('E12').Value = '  +2.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.210'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7113'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.36'
$ws.Range('E15').Value = '  +1.41%  '
$ws.Range('D16').Value = '29.346.77'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008291'
$ws.Range('E17').Value = '  +6.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.004'
$ws.Range('E18').Value = '  +2.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.52'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.127.07'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.22'
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1610'
$ws.Range('E25').Value = '  +1.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.62'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.015'
$ws.Range('E27').Value = '  +1.36%  '
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('E29').Value = '  +1.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.403'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.343'
$ws.Range('E31').Value = '  +5.35%  '
$ws.Range('E32').Value = '  -3.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05259'
$ws.Range('E33').Value = '  +1.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.932'
$ws.Range('E34').Value = '  +1.31%  '
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7430'
$ws.Range('E36').Value = '  +2.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.685'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01868'
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.721'
$ws.Range('E39').Value = '  +1.32%  '
$ws.Range('D40').Value = '1.167.93'
$ws.Range('E40').Value = '  +1.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.381'
$ws.Range('E41').Value = '  +4.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.13'
$ws.Range('E42').Value = '  +1.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8879'
$ws.Range('E43').Value = '  -1.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '106.76'
$ws.Range('E44').Value = '  +4.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9996'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = '2.022.59'
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.816'
$ws.Range('E47').Value = '  +3.07%  '
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000121'
$ws.Range('E49').Value = '  +1.41%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.390'
$ws.Range('E50').Value = '  +1.13%  '
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4303'
$ws.Range('E51').Value = '  +1.42%  '
